$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.379.75'
$ws.Range("E2").Value = '  +1.24%  '

$ws.Range("D3").Value = '2.031.59'
$ws.Range("E3").Value = '  +0.74%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.05%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.29'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.380'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0785'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.103'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.87%  '

$ws.Range("D12").Value = '2.332.39'
$ws.Range("E12").Value = '  +0.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.80%  '

$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.21'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.21%  '

$ws.Range("D17").Value = '2.008.01'
$ws.Range("E17").Value = '  -0.53%  '

$ws.Range("D18").Value = '37.269.22'
$ws.Range("E18").Value = '  +1.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.28%  '

$ws.Range("E21").Value = '  -0.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("E25").Value = '  -0.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.29'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.64%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.14'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.02%  '

$ws.Range("E28").Value = '  +5.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.75'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("E30").Value = '  -1.12%  '

$ws.Range("E31").Value = '  +0.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.48'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.09%  '

$ws.Range("E33").Value = '  -0.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.02'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +10.66%  '

$ws.Range("E35").Value = '  +0.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.52%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.23'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.06%  '

$ws.Range("E38").Value = '  -0.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.66'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.44%  '

$ws.Range("D40").Value = '1.469.82'
$ws.Range("E40").Value = '  -1.13%  '

$ws.Range("E41").Value = '  -2.25%  '

$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0925'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.51%  '

$ws.Range("B43").Value = 'FTXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +14.92%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '94.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.55%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.31'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.97%  '

$ws.Range("E47").Value = '  -2.43%  '

$ws.Range("E48").Value = '  +0.83%  '

$ws.Range("E49").Value = '  -2.41%  '

$ws.Range("E50").Value = '  +1.17%  '

$ws.Range("D51").Value = '2.220.24'
$ws.Range("E51").Value = '  +0.62%  '
